# Insert one new data row at row 102 (shifting the existing rows 102-191
# down to 103-192) and populate it with a new daily price record for
# Ciboulette at Vega Modelo de Temuco, matching the rest of the table's
# layout (dimension grows from A1:R191 to A1:R192).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(102).Insert()

$ws.Cells.Item(102, 1).Value  = 10
$ws.Cells.Item(102, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(102, 3).Value  = "La Araucanía"
$ws.Cells.Item(102, 4).Value  = 44566
$ws.Cells.Item(102, 5).Value  = 9
$ws.Cells.Item(102, 6).Value  = 100112039
$ws.Cells.Item(102, 7).Value  = "Ciboulette"
$ws.Cells.Item(102, 8).Value  = "Sin especificar"
$ws.Cells.Item(102, 9).Value  = "Primera"
$ws.Cells.Item(102, 10).Value = 20
$ws.Cells.Item(102, 11).Value = 5000
$ws.Cells.Item(102, 12).Value = 5000
$ws.Cells.Item(102, 13).Value = 5000
$ws.Cells.Item(102, 14).Value = "$/docena de atados"
$ws.Cells.Item(102, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(102, 16).Value = 1667
$ws.Cells.Item(102, 17).Value = 3
$ws.Cells.Item(102, 18).Value = "Hortaliza"
